# Generate Report for Handoff
# A later handoff run advanced the "11ee4a49-..." (.md) localization entry's
# handoff timestamps, which ripples into:
#   - zh-cn sheet: row for 11ee4a49-...-d51267285a0d, Latest Handoff Datetime (E5)
#   - de-de sheet: row for 11ee4a49-...-d51267285a0d, Latest Handoff Datetime (E5)
#   - Overview sheet: last row's (521bf458-...) Latest Handoff Date (D7)

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E5").Value = "2016-03-23 10:43:12"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E5").Value = "2016-03-23 10:43:16"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D7").Value = "2016-03-23 10:43:16"
